# Refactor currency conversion sheet: split the old single "foreign_amount"
# column into explicit "source_amount" / "target_amount" columns, and add a
# new "target_fees" column next to the existing "source_fees" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("currency_conversions")

# Make room for two new columns (target_amount, target_fees) by inserting
# them where target_currency currently is (columns E:F), shifting the
# existing target_currency/comment columns to the right. Using a native
# column Insert (rather than just writing into new cells) preserves the
# original header cell formatting for the shifted cells.
$ws.Columns.Item(5).Insert(-4161)
$ws.Columns.Item(5).Insert(-4161)

# Rename the old "foreign_amount" header to "source_amount".
$ws.Range("B1").Value = "source_amount"

# Fill in the two newly inserted headers.
$ws.Range("E1").Value = "target_amount"
$ws.Range("F1").Value = "target_fees"

# The column insert above pushed the sheet's trailing (empty) columns past
# the real worksheet boundary; delete the same number of now-superfluous
# trailing columns to restore the normal column range.
$ws.Columns.Item(16383).Delete(-4159)
$ws.Columns.Item(16383).Delete(-4159)

# The currency_conversions sheet is now the one the user was last working
# on, so make it the active/selected sheet (previously money_transfers was
# active).
$ws.Activate()
